$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.727.81'
$ws.Range("E2").Value = '  -5.24%  '

$ws.Range("D3").Value = '3.380.20'
$ws.Range("E3").Value = '  -4.71%  '

$ws.Range("D5").Value = "'555.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.49%  '

$ws.Range("D6").Value = "'170.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -9.58%  '

$ws.Range("D7").Value = "'0.612"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.02%  '

$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").Value = "'0.612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.24%  '

$ws.Range("D10").Value = "'0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.22%  '

$ws.Range("D11").Value = "'55.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.56%  '

$ws.Range("D12").Value = "'0.0000265"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.39%  '

$ws.Range("D13").Value = "'8.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.88%  '

$ws.Range("D14").Value = '3.948.54'
$ws.Range("E14").Value = '  -4.00%  '

$ws.Range("D15").Value = '3.428.87'
$ws.Range("E15").Value = '  -3.34%  '

$ws.Range("E16").Value = '  -2.62%  '

$ws.Range("D17").Value = "'17.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.94%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '63.821.49'
$ws.Range("E18").Value = '  -5.15%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = "'11.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.49%  '

$ws.Range("D20").Value = "'0.979"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.09%  '

$ws.Range("D21").Value = "'406.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.44%  '

$ws.Range("D22").Value = "'4.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.28%  '

$ws.Range("D23").Value = "'4.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.91%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = "'82.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.44%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = "'13.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.50%  '

$ws.Range("D26").Value = "'10.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.61%  '

$ws.Range("D27").Value = "'2.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.95%  '

$ws.Range("D28").Value = "'8.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.78%  '

$ws.Range("D29").Value = "'29.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.97%  '

$ws.Range("D30").Value = "'6.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.94%  '

$ws.Range("D31").Value = "'585.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.30%  '

$ws.Range("D32").Value = "'11.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.13%  '

$ws.Range("D33").Value = "'0.106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.52%  '

$ws.Range("D34").Value = "'58.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.92%  '

$ws.Range("D35").Value = "'0.149"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.40%  '

$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("D37").Value = "'35.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.80%  '

$ws.Range("D38").Value = '3.177.31'
$ws.Range("E38").Value = '  +0.72%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = "'3.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.89%  '

$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = "'0.368"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.11%  '

$ws.Range("D41").Value = '0.0₃0722'
$ws.Range("E41").Value = '  -11.26%  '

$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.07%  '

$ws.Range("D43").Value = "'2.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.12%  '

$ws.Range("D44").Value = "'2.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.10%  '

$ws.Range("D45").Value = "'3.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.22%  '

$ws.Range("D46").Value = "'0.0402"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.99%  '

$ws.Range("D47").Value = "'2.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.64%  '

$ws.Range("D48").Value = "'0.128"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.50%  '

$ws.Range("D49").Value = "'134.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.57%  '

$ws.Range("D50").Value = "'8.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.29%  '

$ws.Range("D51").Value = "'2.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.31%  '
